# Aplikacja Sports Joiner.pptx - timeline / schedule table update (slide 9)
#
# Changes applied:
#  1. Table "Tabela 8" (slide 9, shape 4) row "Najwazniejsze zadania":
#     every bulleted paragraph switches its bullet from the Arial "*"
#     character bullet to a Wingdings "0xD8" ("O-slash") character bullet.
#  2. Table row "Zaleznosci" (last row) gets new text in the four data
#     columns (Planowanie / Kodowanie / Testowanie / Kampania marketingowa).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$tbl = $s.Shapes.Item(4).Table

# --- 1. Bullet restyle: row 3 ("Najwazniejsze zadania"), columns 2-5 ---
for ($col = 2; $col -le 5; $col++) {
    $cell = $tbl.Cell(3, $col)
    $tr = $cell.Shape.TextFrame.TextRange
    $paraCount = $tr.Paragraphs().Count
    for ($i = 1; $i -le $paraCount; $i++) {
        $para = $tr.Paragraphs().Item($i)
        $bullet = $para.ParagraphFormat.Bullet
        # Font must be set before Character so the serialized XML keeps
        # <a:buFont> before <a:buChar> (schema element order).
        $bullet.Font.Name = "Wingdings"
        $bullet.Character = 216
    }
}

# --- 2. Fill in the "Zaleznosci" row (row 4) ---
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "-"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = "-"
$tbl.Cell(4, 4).Shape.TextFrame.TextRange.Text = "Oczekiwanie na gotową aplikację"
$tbl.Cell(4, 5).Shape.TextFrame.TextRange.Text = "Oczekiwanie na zakończone powodzeniem testy"
